# Apply the text replacements described by the diff.
$d = $word.ActiveDocument

$replacements = @(
    @("2025-04-30 Wednesday", "2025-05-01 Thursday"),
    @("164÷7=", "152÷2="),
    @("890÷6=", "342÷9="),
    @("492÷6=", "775÷2="),
    @("565÷9=", "928÷7="),
    @("423÷6=", "874÷8="),
    @("340÷4=", "713÷9="),
    @("922÷7=", "712÷3="),
    @("996÷4=", "542÷6="),
    @("100÷9=", "108÷2="),
    @("225÷4=", "896÷9="),
    @("188÷8=", "334÷3="),
    @("394÷5=", "587÷8="),
    @("485÷9=", "540÷7="),
    @("651÷4=", "721÷6="),
    @("275÷3=", "742÷4="),
    @("465÷6=", "180÷5="),
    @("234÷2=", "587÷4="),
    @("571÷9=", "342÷4="),
    @("624÷9=", "478÷9="),
    @("731÷2=", "248÷4="),
    @("641÷5=", "403÷7="),
    @("458÷2=", "766÷2="),
    @("793÷5=", "550÷9="),
    @("332÷6=", "298÷6="),
    @("363÷2=", "534÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
